$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 3979.4102
$ws.Range("J17").Value = 1433.6884
$ws.Range("L17").Value = 4301.0652
$ws.Range("N17").Value = -4637.0652
$ws.Range("H62").Value = 2908.75
$ws.Range("I62").Value = 2301
$ws.Range("J62").Value = 3921.6667
$ws.Range("K62").Value = 2301
$ws.Range("L62").Value = 3921.6667
$ws.Range("M62").Value = -1677
$ws.Range("N62").Value = -5169.6667
$ws.Range("H65").Value = 2908.75
$ws.Range("I65").Value = 2301
$ws.Range("J65").Value = 3921.6667
$ws.Range("K65").Value = 11505
$ws.Range("L65").Value = 19608.3335
$ws.Range("M65").Value = -8385
$ws.Range("N65").Value = -25848.3335
$ws.Range("H74").Value = 4082.5
$ws.Range("I74").Value = 3795
$ws.Range("J74").Value = 4178.3335
$ws.Range("K74").Value = 3795
$ws.Range("L74").Value = 4178.3335
$ws.Range("M74").Value = -2859
$ws.Range("N74").Value = -6050.3335
$ws.Range("H77").Value = 4082.5
$ws.Range("I77").Value = 3795
$ws.Range("J77").Value = 4178.3335
$ws.Range("K77").Value = 18975
$ws.Range("L77").Value = 20891.6675
$ws.Range("M77").Value = -14295
$ws.Range("N77").Value = -30251.6675
$ws.Range("H107").Value = 679.0526
$ws.Range("I107").Value = 725.4666999999999
$ws.Range("K107").Value = 725.4666999999999
$ws.Range("M107").Value = 1194.5333
$ws.Range("H132").Value = 3336512.8
$ws.Range("I132").Value = 3451386.5
$ws.Range("J132").Value = 5175
$ws.Range("K132").Value = 10354159.5
$ws.Range("L132").Value = 15525
$ws.Range("M132").Value = -10351629.5
$ws.Range("N132").Value = -20585
$ws.Range("H137").Value = 4170945
$ws.Range("I137").Value = 7697961.5
$ws.Range("J137").Value = 2652.6365
$ws.Range("K137").Value = 23093884.5
$ws.Range("L137").Value = 7957.9095
$ws.Range("M137").Value = -23091334.5
$ws.Range("N137").Value = -13057.9095

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 2809.5112
$ws.Range("I61").Value = 1609.5238
$ws.Range("J61").Value = 3859.5
$ws.Range("K61").Value = 1609.5238
$ws.Range("L61").Value = 3859.5
$ws.Range("M61").Value = -1397.5238
$ws.Range("N61").Value = -4283.5
$ws.Range("H74").Value = 519.9474
$ws.Range("I74").Value = 519.9474
$ws.Range("K74").Value = 519.9474
$ws.Range("M74").Value = 354.0526
$ws.Range("H77").Value = 519.9474
$ws.Range("I77").Value = 519.9474
$ws.Range("K77").Value = 2599.737
$ws.Range("M77").Value = 1768.263
$ws.Range("H110").Value = 2450
$ws.Range("I110").Value = 700
$ws.Range("J110").Value = 2609.0908
$ws.Range("K110").Value = 700
$ws.Range("L110").Value = 2609.0908
$ws.Range("M110").Value = 1345
$ws.Range("N110").Value = -6699.0908
$ws.Range("H124").Value = 27000
$ws.Range("J124").Value = 27000
$ws.Range("L124").Value = 27000
$ws.Range("N124").Value = -36820
$ws.Range("H125").Value = 28722.857
$ws.Range("J125").Value = 28722.857
$ws.Range("L125").Value = 28722.857
$ws.Range("N125").Value = -38562.857
$ws.Range("H136").Value = 2809.5112
$ws.Range("I136").Value = 1609.5238
$ws.Range("J136").Value = 3859.5
$ws.Range("K136").Value = 4828.5714
$ws.Range("L136").Value = 11578.5
$ws.Range("M136").Value = -2278.5714
$ws.Range("N136").Value = -16678.5

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1984.5294
$ws.Range("I20").Value = 1000
$ws.Range("J20").Value = 2521.5454
$ws.Range("K20").Value = 1000
$ws.Range("L20").Value = 2521.5454
$ws.Range("M20").Value = -753
$ws.Range("N20").Value = -3015.5454
$ws.Range("H86").Value = 22735.875
$ws.Range("I86").Value = 1129.4286
$ws.Range("K86").Value = 1129.4286
$ws.Range("M86").Value = -6.42859999999996
$ws.Range("H89").Value = 22735.875
$ws.Range("I89").Value = 1129.4286
$ws.Range("K89").Value = 5647.143
$ws.Range("M89").Value = -31.14300000000003
$ws.Range("H94").Value = 528.46344
$ws.Range("I94").Value = 501.27026
$ws.Range("J94").Value = 780
$ws.Range("K94").Value = 501.27026
$ws.Range("L94").Value = 780
$ws.Range("M94").Value = -50.27026000000001
$ws.Range("N94").Value = -1682
$ws.Range("H105").Value = 1697.122
$ws.Range("I105").Value = 1468.8235
$ws.Range("J105").Value = 1858.8334
$ws.Range("K105").Value = 1468.8235
$ws.Range("L105").Value = 1858.8334
$ws.Range("M105").Value = 278.1765
$ws.Range("N105").Value = -5352.8334

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2443578
$ws.Range("I31").Value = 3033687.8
$ws.Range("J31").Value = 9375
$ws.Range("K31").Value = 3033687.8
$ws.Range("L31").Value = 9375
$ws.Range("M31").Value = -3033392.8
$ws.Range("N31").Value = -9965
$ws.Range("H34").Value = 2443578
$ws.Range("I34").Value = 3033687.8
$ws.Range("J34").Value = 9375
$ws.Range("K34").Value = 3033687.8
$ws.Range("L34").Value = 9375
$ws.Range("M34").Value = -3033485.8
$ws.Range("N34").Value = -9779

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H88").Value = 4669.3335
$ws.Range("J88").Value = 4669.3335
$ws.Range("L88").Value = 14008.0005
$ws.Range("N88").Value = -14864.0005
$ws.Range("H91").Value = 4669.3335
$ws.Range("J91").Value = 4669.3335
$ws.Range("L91").Value = 14008.0005
$ws.Range("N91").Value = -16972.0005
$ws.Range("H131").Value = 1267.1
$ws.Range("I131").Value = 6500
$ws.Range("J131").Value = 1049.0625
$ws.Range("K131").Value = 19500
$ws.Range("L131").Value = 3147.1875
$ws.Range("M131").Value = -14460
$ws.Range("N131").Value = -13227.1875

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 80109.46000000001
$ws.Range("I102").Value = 2542.4
$ws.Range("K102").Value = 2542.4
$ws.Range("M102").Value = -920.4000000000001

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 3398.5715
$ws.Range("J7").Value = 3298.3333
$ws.Range("L7").Value = 3298.3333
$ws.Range("N7").Value = -3522.3333
$ws.Range("H40").Value = 3435.8
$ws.Range("I40").Value = 1800
$ws.Range("J40").Value = 3844.75
$ws.Range("K40").Value = 1800
$ws.Range("L40").Value = 3844.75
$ws.Range("M40").Value = -1664
$ws.Range("N40").Value = -4116.75
$ws.Range("H122").Value = 2845.8235
$ws.Range("I122").Value = 2418.182
$ws.Range("J122").Value = 3629.8333
$ws.Range("K122").Value = 7254.545999999999
$ws.Range("L122").Value = 10889.4999
$ws.Range("M122").Value = -4804.545999999999
$ws.Range("N122").Value = -15789.4999
$ws.Range("H126").Value = 3398.5715
$ws.Range("J126").Value = 3298.3333
$ws.Range("L126").Value = 9894.999899999999
$ws.Range("N126").Value = -14834.9999
$ws.Range("H127").Value = 32000
$ws.Range("J127").Value = 32000
$ws.Range("L127").Value = 32000
$ws.Range("N127").Value = -41920
$ws.Range("H132").Value = 2776.8108
$ws.Range("I132").Value = 1763.2778
$ws.Range("J132").Value = 3737
$ws.Range("K132").Value = 5289.8334
$ws.Range("L132").Value = 11211
$ws.Range("M132").Value = -2759.8334
$ws.Range("N132").Value = -16271

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 710.44446
$ws.Range("I81").Value = 649
$ws.Range("J81").Value = 833.3333
$ws.Range("K81").Value = 1298
$ws.Range("L81").Value = 1666.6666
$ws.Range("M81").Value = -237
$ws.Range("N81").Value = -3788.6666
$ws.Range("H84").Value = 710.44446
$ws.Range("I84").Value = 649
$ws.Range("J84").Value = 833.3333
$ws.Range("K84").Value = 6490
$ws.Range("L84").Value = 8333.333000000001
$ws.Range("M84").Value = -1186
$ws.Range("N84").Value = -18941.333
$ws.Range("H113").Value = 1064.24
$ws.Range("I113").Value = 444.2143
$ws.Range("J113").Value = 1853.3636
$ws.Range("K113").Value = 1332.6429
$ws.Range("L113").Value = 5560.0908
$ws.Range("M113").Value = 837.3571000000002
$ws.Range("N113").Value = -9900.0908
$ws.Range("H122").Value = 528208.4399999999
$ws.Range("I122").Value = 715686.8
$ws.Range("J122").Value = 3269
$ws.Range("K122").Value = 2147060.4
$ws.Range("L122").Value = 9807
$ws.Range("M122").Value = -2144610.4
$ws.Range("N122").Value = -14707
$ws.Range("H132").Value = 190189.33
$ws.Range("I132").Value = 252655
$ws.Range("K132").Value = 757965
$ws.Range("M132").Value = -755435
